$wb = $excel.ActiveWorkbook

# --- normality sheet ---
$ws = $wb.Worksheets.Item("normality")
$ws.Range("C4").Value = 0.9746
$ws.Range("D4").Value = 0.5979
$ws.Range("C7").Value = 0.8815
$ws.Range("D7").Value = 0.0015
$ws.Range("C10").Value = 0.8744
$ws.Range("D10").Value = 0.001
$ws.Range("C13").Value = 0.9143
$ws.Range("D13").Value = 0.0112
$ws.Range("C16").Value = 0.9061
$ws.Range("D16").Value = 0.0067

# --- equal_var sheet ---
$ws = $wb.Worksheets.Item("equal_var")
$ws.Range("C3").Value = 3.3268
$ws.Range("D3").Value = 0.0398
$ws.Range("C4").Value = 0.4477
$ws.Range("D4").Value = 0.6403
$ws.Range("C5").Value = 1.2232
$ws.Range("D5").Value = 0.2985
$ws.Range("C6").Value = 0.2681
$ws.Range("D6").Value = 0.7654
$ws.Range("C7").Value = 2.9704
$ws.Range("D7").Value = 0.0557

# --- welch_anova sheet ---
$ws = $wb.Worksheets.Item("welch_anova")
$ws.Range("E3").Value = 62.2065
$ws.Range("F3").Value = 5.6844
$ws.Range("G3").Value = 0.0054
$ws.Range("H3").Value = 0.0708

$ws.Range("E4").Value = 53.177
$ws.Range("F4").Value = 0.4525
$ws.Range("G4").Value = 0.6385
$ws.Range("H4").Value = 0.0089

$ws.Range("E5").Value = 61.4906
$ws.Range("F5").Value = 5.0987
$ws.Range("G5").Value = 0.0089
$ws.Range("H5").Value = 0.06270000000000001

$ws.Range("E6").Value = 53.5269
$ws.Range("F6").Value = 2.89
$ws.Range("G6").Value = 0.0643
$ws.Range("H6").Value = 0.053

$ws.Range("E7").Value = 65.328
$ws.Range("F7").Value = 5.7742
$ws.Range("G7").Value = 0.0049
$ws.Range("H7").Value = 0.0762

# --- pairwise_gameshowell sheet ---
$ws = $wb.Worksheets.Item("pairwise_gameshowell")

$ws.Range("E3").Value = 828.1365
$ws.Range("G3").Value = -56.7983
$ws.Range("H3").Value = 61.4999
$ws.Range("I3").Value = -0.9236
$ws.Range("J3").Value = 77.575
$ws.Range("K3").Value = 0.6271
$ws.Range("L3").Value = -0.2034

$ws.Range("E4").Value = 828.1365
$ws.Range("G4").Value = 135.2281
$ws.Range("H4").Value = 59.7676
$ws.Range("I4").Value = 2.2626
$ws.Range("J4").Value = 53.1676
$ws.Range("K4").Value = 0.0701
$ws.Range("L4").Value = 0.6104000000000001

$ws.Range("E6").Value = 0.2588
$ws.Range("G6").Value = 0.0442
$ws.Range("H6").Value = 0.0501
$ws.Range("I6").Value = 0.8833
$ws.Range("J6").Value = 66.36490000000001
$ws.Range("K6").Value = 0.6526999999999999
$ws.Range("L6").Value = 0.1945

$ws.Range("E7").Value = 0.2588
$ws.Range("G7").Value = 0.0091
$ws.Range("H7").Value = 0.0613
$ws.Range("I7").Value = 0.1481
$ws.Range("J7").Value = 47.3108
$ws.Range("K7").Value = 0.988
$ws.Range("L7").Value = 0.04

$ws.Range("E9").Value = 526.6651000000001
$ws.Range("G9").Value = 48.6388
$ws.Range("H9").Value = 47.2102
$ws.Range("I9").Value = 1.0303
$ws.Range("J9").Value = 69.5658
$ws.Range("K9").Value = 0.5604
$ws.Range("L9").Value = 0.2269

$ws.Range("E10").Value = 526.6651000000001
$ws.Range("G10").Value = 141.822
$ws.Range("H10").Value = 46.7897
$ws.Range("I10").Value = 3.0311
$ws.Range("J10").Value = 53.9788
$ws.Range("K10").Value = 0.0103
$ws.Range("L10").Value = 0.8178

$ws.Range("E12").Value = 3.7941
$ws.Range("G12").Value = -0.3859
$ws.Range("H12").Value = 0.3308
$ws.Range("I12").Value = -1.1664
$ws.Range("J12").Value = 77.71850000000001
$ws.Range("K12").Value = 0.4767
$ws.Range("L12").Value = -0.2569

$ws.Range("E13").Value = 3.7941
$ws.Range("G13").Value = -1.0241
$ws.Range("H13").Value = 0.425
$ws.Range("I13").Value = -2.4096
$ws.Range("J13").Value = 39.3487
$ws.Range("K13").Value = 0.053
$ws.Range("L13").Value = -0.6501

$ws.Range("E15").Value = 279.5242
$ws.Range("G15").Value = -70.46729999999999
$ws.Range("H15").Value = 40.1617
$ws.Range("I15").Value = -1.7546
$ws.Range("J15").Value = 80.2805
$ws.Range("K15").Value = 0.1916
$ws.Range("L15").Value = -0.3864

$ws.Range("E16").Value = 279.5242
$ws.Range("G16").Value = 53.8918
$ws.Range("H16").Value = 35.3603
$ws.Range("I16").Value = 1.5241
$ws.Range("J16").Value = 53.9999
$ws.Range("K16").Value = 0.2878
$ws.Range("L16").Value = 0.4112
